# Update metric values in Sheet1 (B2:I26) with new results.
# All 25 data rows (2-26) share the same value per column, per the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

$ws.Range("B2:B26").Value = 0.353672031788087
$ws.Range("C2:C26").Value = 0.4975558783676419
$ws.Range("D2:D26").Value = 0.4947511081832262
$ws.Range("E2:E26").Value = 0.496641305776536
$ws.Range("F2:F26").Value = 0.7152945399284363
$ws.Range("G2:G26").Value = 0.09416297078132629
$ws.Range("H2:H26").Value = 1.526041030883789
$ws.Range("I2:I26").Value = 0.7679885029792786
